$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting the existing "Link" column
# (and its data) from B to C.
$ws.Columns("B").Insert()

# New "Title" column header + values.
$titles = @(
    "Title",
    "X-COM - Wikipedia, la enciclopedia libre",
    "Ahorra un 80 % en XCOM: Enemy Unknown en Steam",
    "XCOM 2",
    "No Title",
    "No Title",
    "No Title",
    "No Title",
    "No Title",
    "Comprar XCOM®: Enemy Unknown",
    "XCOM 2 Collection - Apps en Google Play",
    "Xcom Ps3 | MercadoLibre",
    "Xcom",
    "Todos los juegos de XCOM - Saga completa",
    "XCOM 2"
)

for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}

# Row 15's link URL changed (new srsltid tracking param).
$ws.Range("C15").Value = "https://store.2k.com/es/game/buy-xcom-2-pc?srsltid=AfmBOoqcvQrXcFeuNrafqptyWYzHAs-rjukDGJTw5LRqoiLUbuywdmj7"
